# Auto-generated edit script: update cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "37.917.25"
$ws.Cells.Item(2, 5).Value = "  -0.89%  "
$ws.Cells.Item(3, 4).Value = "2.046.92"
$ws.Cells.Item(3, 5).Value = "  -0.49%  "
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
$ws.Cells.Item(5, 4).Value = "228.73"
$ws.Cells.Item(5, 5).Value = "  -0.21%  "
$ws.Cells.Item(6, 5).Value = "  -1.07%  "
$ws.Cells.Item(7, 4).Value = "61.52"
$ws.Cells.Item(7, 5).Value = "  +0.82%  "
$ws.Cells.Item(8, 5).Value = "  +0.03%  "
$ws.Cells.Item(9, 4).Value = "0.378"
$ws.Cells.Item(9, 5).Value = "  -2.31%  "
$ws.Cells.Item(10, 4).Value = "0.0819"
$ws.Cells.Item(10, 5).Value = "  -1.08%  "
$ws.Cells.Item(12, 2).Value = "Chainlink"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(12, 4).Value = "14.68"
$ws.Cells.Item(12, 5).Value = "  -1.16%  "
$ws.Cells.Item(13, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(13, 4).Value = "2.350.43"
$ws.Cells.Item(13, 5).Value = "  -0.49%  "
$ws.Cells.Item(14, 4).Value = "21.15"
$ws.Cells.Item(14, 5).Value = "  -0.05%  "
$ws.Cells.Item(15, 4).Value = "0.777"
$ws.Cells.Item(15, 5).Value = "  +1.62%  "
$ws.Cells.Item(16, 5).Value = "  -2.39%  "
$ws.Cells.Item(17, 4).Value = "2.074.42"
$ws.Cells.Item(17, 5).Value = "  +0.69%  "
$ws.Cells.Item(18, 4).Value = "37.867.70"
$ws.Cells.Item(18, 5).Value = "  -0.89%  "
$ws.Cells.Item(19, 4).Value = "69.65"
$ws.Cells.Item(20, 5).Value = "  -5.56%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0824"
$ws.Cells.Item(21, 5).Value = "  -1.48%  "
$ws.Cells.Item(22, 4).Value = "224.09"
$ws.Cells.Item(22, 5).Value = "  -0.57%  "
$ws.Cells.Item(23, 5).Value = "  +0.03%  "
$ws.Cells.Item(24, 4).Value = "2.43"
$ws.Cells.Item(24, 5).Value = "  -0.11%  "
$ws.Cells.Item(25, 5).Value = "  +2.65%  "
$ws.Cells.Item(26, 4).Value = "168.26"
$ws.Cells.Item(26, 5).Value = "  +1.15%  "
$ws.Cells.Item(27, 4).Value = "9.36"
$ws.Cells.Item(27, 5).Value = "  +0.76%  "
$ws.Cells.Item(28, 4).Value = "0.129"
$ws.Cells.Item(28, 5).Value = "  -2.18%  "
$ws.Cells.Item(29, 4).Value = "18.86"
$ws.Cells.Item(29, 5).Value = "  -0.71%  "
$ws.Cells.Item(30, 5).Value = "  -1.20%  "
$ws.Cells.Item(31, 4).Value = "0.121"
$ws.Cells.Item(31, 5).Value = "  +0.20%  "
$ws.Cells.Item(32, 5).Value = "  +8.16%  "
$ws.Cells.Item(33, 4).Value = "'4.40"
$ws.Cells.Item(33, 5).Value = "  -2.43%  "
$ws.Cells.Item(34, 5).Value = "  -1.28%  "
$ws.Cells.Item(35, 5).Value = "  +0.47%  "
$ws.Cells.Item(36, 4).Value = "6.64"
$ws.Cells.Item(36, 5).Value = "  +3.89%  "
$ws.Cells.Item(37, 5).Value = "  +2.02%  "
$ws.Cells.Item(38, 4).Value = "3.45"
$ws.Cells.Item(38, 5).Value = "  +5.17%  "
$ws.Cells.Item(39, 5).Value = "  +0.04%  "
$ws.Cells.Item(40, 4).Value = "18.18"
$ws.Cells.Item(40, 5).Value = "  +7.37%  "
$ws.Cells.Item(41, 4).Value = "1.544.96"
$ws.Cells.Item(41, 5).Value = "  +0.99%  "
$ws.Cells.Item(42, 5).Value = "  +0.26%  "
$ws.Cells.Item(43, 4).Value = "96.58"
$ws.Cells.Item(43, 5).Value = "  -1.31%  "
$ws.Cells.Item(44, 5).Value = "  -1.78%  "
$ws.Cells.Item(45, 4).Value = "0.0917"
$ws.Cells.Item(45, 5).Value = "  -1.62%  "
$ws.Cells.Item(46, 4).Value = "4.13"
$ws.Cells.Item(46, 5).Value = "  +2.99%  "
$ws.Cells.Item(47, 5).Value = "  -1.45%  "
$ws.Cells.Item(48, 5).Value = "  -0.35%  "
$ws.Cells.Item(49, 5).Value = "  -0.31%  "
$ws.Cells.Item(50, 4).Value = "7.06"
$ws.Cells.Item(50, 5).Value = "  +0.08%  "
$ws.Cells.Item(51, 4).Value = "2.239.51"
$ws.Cells.Item(51, 5).Value = "  -0.48%  "
